# Corrects the Makeham mortality law parameters (A, B) used to build the
# life table on the "makeham" sheet. c (the geometric-growth factor of the
# force of mortality) is unchanged at 1.07; A: 0.000002 -> 0.0001 and
# B: 0.000006 -> 0.0003 (both x50). Columns: A=x, B=lx, C=dx, D=qx, E=px,
# F=exo (complete expectation of life). Row r holds age x = r-2, for
# r = 2..128 (x = 0..126, omega = 126).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mkA = 0.0001
$mkB = 0.0003
$mkc = 1.07
$l0 = 100000
$omega = 126

$lnc = $ws.Evaluate("LN($mkc)")

$qx = New-Object 'double[]' ($omega + 1)
$lx = New-Object 'double[]' ($omega + 1)
$dx = New-Object 'double[]' ($omega + 1)
$px = New-Object 'double[]' ($omega + 1)
$ex = New-Object 'double[]' ($omega + 1)

# qx(x) = 1 - exp( -( A + B*c^x*(c-1)/ln(c) ) )  for x = 0 .. omega-1
for ($x = 0; $x -lt $omega; $x++) {
    $formula = "1-EXP(-($mkA+$mkB*$mkc^$x*($mkc-1)/$lnc))"
    $qx[$x] = $ws.Evaluate($formula)
}
$qx[$omega] = 1

# lx(0) = l0 ; lx(x+1) = lx(x) * (1 - qx(x))
$lx[0] = $l0
for ($x = 0; $x -lt $omega; $x++) {
    $lx[$x + 1] = $lx[$x] * (1 - $qx[$x])
}

# dx(x) = lx(x) * qx(x) for x < omega ; dx(omega) = lx(omega)
for ($x = 0; $x -lt $omega; $x++) {
    $dx[$x] = $lx[$x] * $qx[$x]
}
$dx[$omega] = $lx[$omega]

# px(x) = 1 - qx(x) for x < omega ; px(omega) = 0
for ($x = 0; $x -lt $omega; $x++) {
    $px[$x] = 1 - $qx[$x]
}
$px[$omega] = 0

# exo(x) = complete expectation of life = curtate e_x + 0.5
#        = ( sum_{k=x+1}^{omega} lx(k) ) / lx(x) + 0.5 , for x < omega
# exo(omega) = 0.5
$tailSum = 0
for ($x = $omega - 1; $x -ge 0; $x--) {
    $tailSum = $tailSum + $lx[$x + 1]
    $ex[$x] = $tailSum / $lx[$x] + 0.5
}
$ex[$omega] = 0.5

# Write the recomputed table back onto the sheet (rows 2..128, x = 0..126)
for ($x = 0; $x -le $omega; $x++) {
    $r = $x + 2
    $ws.Cells.Item($r, 2).Value = $lx[$x]
    $ws.Cells.Item($r, 3).Value = $dx[$x]
    $ws.Cells.Item($r, 4).Value = $qx[$x]
    $ws.Cells.Item($r, 5).Value = $px[$x]
    $ws.Cells.Item($r, 6).Value = $ex[$x]
}
